$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new data row (row 4) mirroring row 3 but with a new username "paul123"
$ws.Range("A4").Value = "paul123"
$ws.Range("B4").Value = "peter"
$ws.Range("C4").Value = "Dutch"
$ws.Range("D4").Value = "Invalid username or password"

# Update selection to reflect the new active range
$ws.Range("D3:D4").Select()
